$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 5133.7144
$ws.Cells.Item(82, 9).Value = 1460.5
$ws.Cells.Item(82, 10).Value = 10031.333
$ws.Cells.Item(82, 11).Value = 4381.5
$ws.Cells.Item(82, 12).Value = 30093.999
$ws.Cells.Item(82, 13).Value = -3975.5
$ws.Cells.Item(82, 14).Value = -30905.999
$ws.Cells.Item(85, 8).Value = 5133.7144
$ws.Cells.Item(85, 9).Value = 1460.5
$ws.Cells.Item(85, 10).Value = 10031.333
$ws.Cells.Item(85, 11).Value = 4381.5
$ws.Cells.Item(85, 12).Value = 30093.999
$ws.Cells.Item(85, 13).Value = -2977.5
$ws.Cells.Item(85, 14).Value = -32901.999
$ws.Cells.Item(86, 8).Value = 41920
$ws.Cells.Item(86, 9).Value = 27400
$ws.Cells.Item(86, 10).Value = 100000
$ws.Cells.Item(86, 11).Value = 27400
$ws.Cells.Item(86, 12).Value = 100000
$ws.Cells.Item(86, 13).Value = -26277
$ws.Cells.Item(86, 14).Value = -102246
$ws.Cells.Item(87, 8).Value = 35000
$ws.Cells.Item(87, 10).Value = 35000
$ws.Cells.Item(87, 12).Value = 35000
$ws.Cells.Item(87, 14).Value = -37496
$ws.Cells.Item(88, 8).Value = 1588630.4
$ws.Cells.Item(88, 9).Value = 3704137.8
$ws.Cells.Item(88, 10).Value = 2000
$ws.Cells.Item(88, 11).Value = 3704137.8
$ws.Cells.Item(88, 12).Value = 2000
$ws.Cells.Item(88, 13).Value = -3703731.8
$ws.Cells.Item(88, 14).Value = -2812
$ws.Cells.Item(89, 8).Value = 41920
$ws.Cells.Item(89, 9).Value = 27400
$ws.Cells.Item(89, 10).Value = 100000
$ws.Cells.Item(89, 11).Value = 137000
$ws.Cells.Item(89, 12).Value = 500000
$ws.Cells.Item(89, 13).Value = -131384
$ws.Cells.Item(89, 14).Value = -511232
$ws.Cells.Item(90, 8).Value = 35000
$ws.Cells.Item(90, 10).Value = 35000
$ws.Cells.Item(90, 12).Value = 105000
$ws.Cells.Item(90, 14).Value = -117480
$ws.Cells.Item(91, 8).Value = 1588630.4
$ws.Cells.Item(91, 9).Value = 3704137.8
$ws.Cells.Item(91, 10).Value = 2000
$ws.Cells.Item(91, 11).Value = 3704137.8
$ws.Cells.Item(91, 12).Value = 2000
$ws.Cells.Item(91, 13).Value = -3702733.8
$ws.Cells.Item(91, 14).Value = -4808
$ws.Cells.Item(132, 8).Value = 2011.0333
$ws.Cells.Item(132, 9).Value = 1358.3334
$ws.Cells.Item(132, 11).Value = 4075.0002
$ws.Cells.Item(132, 13).Value = -1545.0002
$ws.Cells.Item(135, 8).Value = 4456.28
$ws.Cells.Item(135, 9).Value = 5901.0557
$ws.Cells.Item(135, 11).Value = 53109.5013
$ws.Cells.Item(135, 13).Value = -50574.5013
$ws.Cells.Item(137, 8).Value = 1139.4546
$ws.Cells.Item(137, 9).Value = 667.5
$ws.Cells.Item(137, 10).Value = 1244.3334
$ws.Cells.Item(137, 11).Value = 2002.5
$ws.Cells.Item(137, 12).Value = 3733.0002
$ws.Cells.Item(137, 13).Value = 547.5
$ws.Cells.Item(137, 14).Value = -8833.0002
$ws.Cells.Item(141, 8).Value = 4403.7095
$ws.Cells.Item(141, 9).Value = 2475.2778
$ws.Cells.Item(141, 11).Value = 7425.8334
$ws.Cells.Item(141, 13).Value = -2245.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11220.412
$ws.Cells.Item(32, 9).Value = 3827.468
$ws.Cells.Item(32, 10).Value = 27766.523
$ws.Cells.Item(32, 11).Value = 3827.468
$ws.Cells.Item(32, 12).Value = 27766.523
$ws.Cells.Item(32, 13).Value = -3540.468
$ws.Cells.Item(32, 14).Value = -28340.523
$ws.Cells.Item(61, 8).Value = 1304.7307
$ws.Cells.Item(61, 9).Value = 517.26666
$ws.Cells.Item(61, 10).Value = 2378.5454
$ws.Cells.Item(61, 11).Value = 517.26666
$ws.Cells.Item(61, 12).Value = 2378.5454
$ws.Cells.Item(61, 13).Value = -305.26666
$ws.Cells.Item(61, 14).Value = -2802.5454
$ws.Cells.Item(132, 8).Value = 1838.38
$ws.Cells.Item(132, 9).Value = 938.43335
$ws.Cells.Item(132, 10).Value = 3188.3
$ws.Cells.Item(132, 11).Value = 2815.30005
$ws.Cells.Item(132, 12).Value = 9564.900000000001
$ws.Cells.Item(132, 13).Value = -285.3000499999998
$ws.Cells.Item(132, 14).Value = -14624.9
$ws.Cells.Item(136, 8).Value = 1304.7307
$ws.Cells.Item(136, 9).Value = 517.26666
$ws.Cells.Item(136, 10).Value = 2378.5454
$ws.Cells.Item(136, 11).Value = 1551.79998
$ws.Cells.Item(136, 12).Value = 7135.6362
$ws.Cells.Item(136, 13).Value = 998.20002
$ws.Cells.Item(136, 14).Value = -12235.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1673.3871
$ws.Cells.Item(58, 9).Value = 1373.1875
$ws.Cells.Item(58, 10).Value = 1993.6
$ws.Cells.Item(58, 11).Value = 1373.1875
$ws.Cells.Item(58, 12).Value = 1993.6
$ws.Cells.Item(58, 13).Value = -1170.1875
$ws.Cells.Item(58, 14).Value = -2399.6
$ws.Cells.Item(132, 8).Value = 4800.75
$ws.Cells.Item(132, 9).Value = 4512.3335
$ws.Cells.Item(132, 10).Value = 5666
$ws.Cells.Item(132, 11).Value = 13537.0005
$ws.Cells.Item(132, 12).Value = 16998
$ws.Cells.Item(132, 13).Value = -11007.0005
$ws.Cells.Item(132, 14).Value = -22058
$ws.Cells.Item(134, 8).Value = 1134.6182
$ws.Cells.Item(134, 9).Value = 1097.3143
$ws.Cells.Item(134, 10).Value = 1199.9
$ws.Cells.Item(134, 11).Value = 3291.9429
$ws.Cells.Item(134, 12).Value = 3599.7
$ws.Cells.Item(134, 13).Value = -756.9429
$ws.Cells.Item(134, 14).Value = -8669.700000000001
$ws.Cells.Item(136, 8).Value = 1673.3871
$ws.Cells.Item(136, 9).Value = 1373.1875
$ws.Cells.Item(136, 10).Value = 1993.6
$ws.Cells.Item(136, 11).Value = 4119.5625
$ws.Cells.Item(136, 12).Value = 5980.799999999999
$ws.Cells.Item(136, 13).Value = -1569.5625
$ws.Cells.Item(136, 14).Value = -11080.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1471.2805
$ws.Cells.Item(68, 9).Value = 1128.878
$ws.Cells.Item(68, 10).Value = 1813.683
$ws.Cells.Item(68, 11).Value = 3386.634
$ws.Cells.Item(68, 12).Value = 5441.049
$ws.Cells.Item(68, 13).Value = -2575.634
$ws.Cells.Item(68, 14).Value = -7063.049
$ws.Cells.Item(71, 8).Value = 1471.2805
$ws.Cells.Item(71, 9).Value = 1128.878
$ws.Cells.Item(71, 10).Value = 1813.683
$ws.Cells.Item(71, 11).Value = 10159.902
$ws.Cells.Item(71, 12).Value = 16323.147
$ws.Cells.Item(71, 13).Value = -6103.902
$ws.Cells.Item(71, 14).Value = -24435.147
$ws.Cells.Item(107, 8).Value = 563.7941
$ws.Cells.Item(107, 9).Value = 373.75
$ws.Cells.Item(107, 10).Value = 1450.6666
$ws.Cells.Item(107, 11).Value = 1121.25
$ws.Cells.Item(107, 12).Value = 4351.9998
$ws.Cells.Item(107, 13).Value = 798.75
$ws.Cells.Item(107, 14).Value = -8191.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1991.3334
$ws.Cells.Item(122, 9).Value = 2116.2
$ws.Cells.Item(122, 10).Value = 1741.6
$ws.Cells.Item(122, 11).Value = 6348.599999999999
$ws.Cells.Item(122, 12).Value = 5224.799999999999
$ws.Cells.Item(122, 13).Value = -3898.599999999999
$ws.Cells.Item(122, 14).Value = -10124.8
$ws.Cells.Item(132, 8).Value = 4703.162
$ws.Cells.Item(132, 9).Value = 5864.5415
$ws.Cells.Item(132, 10).Value = 2559.077
$ws.Cells.Item(132, 11).Value = 17593.6245
$ws.Cells.Item(132, 12).Value = 7677.231000000001
$ws.Cells.Item(132, 13).Value = -15063.6245
$ws.Cells.Item(132, 14).Value = -12737.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2574.4146
$ws.Cells.Item(132, 9).Value = 1624.48
$ws.Cells.Item(132, 10).Value = 4058.6875
$ws.Cells.Item(132, 11).Value = 4873.440000000001
$ws.Cells.Item(132, 12).Value = 12176.0625
$ws.Cells.Item(132, 13).Value = -2343.440000000001
$ws.Cells.Item(132, 14).Value = -17236.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1046.6578
$ws.Cells.Item(132, 9).Value = 803.4091
$ws.Cells.Item(132, 10).Value = 1381.125
$ws.Cells.Item(132, 11).Value = 2410.2273
$ws.Cells.Item(132, 12).Value = 4143.375
$ws.Cells.Item(132, 13).Value = 119.7727
$ws.Cells.Item(132, 14).Value = -9203.375
